$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Responsibile Parties")
$ws.Range("A4").Value = "Sudo-Kengo"
